$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12614.48 pesos`n✅ 12614.48 pesos = 3.3 = 966.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update numeric cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("O10").Value = 3778
$ws2.Range("N12").Value = 3819
$ws2.Range("O12").Value = 292.7
